$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 288-289, pushing existing rows 288-309 down to 290-311
$ws.Rows("288:289").Insert()

# Copy the date number format (style) from row 290 column D to the new rows' D cells
$ws.Range("D288").NumberFormat = $ws.Range("D290").NumberFormat
$ws.Range("D289").NumberFormat = $ws.Range("D290").NumberFormat

# Row 288 data
$ws.Range("A288").Value = 5
$ws.Range("B288").Value = "Macroferia Regional de Talca"
$ws.Range("C288").Value = "Maule"
$ws.Range("D288").Value = 44783
$ws.Range("E288").Value = 7
$ws.Range("F288").Value = "Fruta"
$ws.Range("G288").Value = 100102
$ws.Range("H288").Value = "Cítricos"
$ws.Range("I288").Value = 100102004
$ws.Range("J288").Value = "Mandarina"
$ws.Range("K288").Value = "Murcott"
$ws.Range("L288").Value = "Segunda"
$ws.Range("M288").Value = 250
$ws.Range("N288").Value = 7000
$ws.Range("O288").Value = 7000
$ws.Range("P288").Value = 7000
$ws.Range("Q288").Value = "$/bandeja 18 kilos"
$ws.Range("R288").Value = "Región de O'Higgins"
$ws.Range("S288").Value = 389
$ws.Range("T288").Value = 18

# Row 289 data
$ws.Range("A289").Value = 5
$ws.Range("B289").Value = "Macroferia Regional de Talca"
$ws.Range("C289").Value = "Maule"
$ws.Range("D289").Value = 44783
$ws.Range("E289").Value = 7
$ws.Range("F289").Value = "Fruta"
$ws.Range("G289").Value = 100102
$ws.Range("H289").Value = "Cítricos"
$ws.Range("I289").Value = 100102004
$ws.Range("J289").Value = "Mandarina"
$ws.Range("K289").Value = "Murcott"
$ws.Range("L289").Value = "Tercera"
$ws.Range("M289").Value = 150
$ws.Range("N289").Value = 5000
$ws.Range("O289").Value = 5000
$ws.Range("P289").Value = 5000
$ws.Range("Q289").Value = "$/bandeja 18 kilos"
$ws.Range("R289").Value = "Región de O'Higgins"
$ws.Range("S289").Value = 278
$ws.Range("T289").Value = 18
